# Added two user stories based off of what we were coming up with in class.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the "enter resources" / cost-tracking story with a new
# "games I own" user story.
$ws.Range("C2").Value = "Add games that I own/enjoy"
$ws.Range("D2").Value = "I can keep track of all the games in my collection"

# Replace the "overview budget" story with a new "communicate with
# others" user story.
$ws.Range("C3").Value = "Communicate with others with the same game"
$ws.Range("D3").Value = "I can find new friends/ get and share (port) information about the game"

$ws.Range("D13").Select()
